$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 30-42 (data no longer present after the update)
$ws.Range("A30:C42").EntireRow.Delete() | Out-Null

# Retraining the 15 minute IGCC netting flow values with the latest data
$data = @(
    @(2, 45923, 10.365, 0.051),
    @(3, 45923.01041666666, 27.322, 0),
    @(4, 45923.02083333334, 29.356, 0),
    @(5, 45923.03125, 22.337, 0),
    @(6, 45923.04166666666, 6.487, 0.144),
    @(7, 45923.05208333334, 15.016, 0.438),
    @(8, 45923.0625, 15.059, 0),
    @(9, 45923.07291666666, 7.71, 0.045),
    @(10, 45923.08333333334, 2.738, 0.635),
    @(11, 45923.09375, 0.464, 4.228),
    @(12, 45923.10416666666, 1.03, 1.049),
    @(13, 45923.11458333334, 3.885, 0.288),
    @(14, 45923.125, 4.566, 0.649),
    @(15, 45923.13541666666, 13.783, 0),
    @(16, 45923.14583333334, 36.378, 0),
    @(17, 45923.15625, 12.307, 0.047),
    @(18, 45923.16666666666, 5.717, 0.019),
    @(19, 45923.17708333334, 5.318, 0.145),
    @(20, 45923.1875, 20.199, 0),
    @(21, 45923.19791666666, 11.381, 0.002),
    @(22, 45923.20833333334, 2.359, 0.113),
    @(23, 45923.21875, 4.97, 0.037),
    @(24, 45923.22916666666, 5.917, 0),
    @(25, 45923.23958333334, 7.482, 0),
    @(26, 45923.25, 0.738, 2.567),
    @(27, 45923.26041666666, 1.362, 0),
    @(28, 45923.27083333334, 2.887, 0),
    @(29, 45923.28125, 0.5659999999999999, 4.542)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
    $ws.Cells.Item($r, 3).Value = $item[3]
}
